# Slide 3 ("slide 2" in the author's 0-based numbering), the content
# placeholder that currently reads a single paragraph "For test" needs to
# become:
#   Paragraph 1: two runs -> "For " + "test"
#   Paragraph 2 (new): one run -> "test"
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# 1) Append a brand-new paragraph containing "test" right after the
#    existing "For test" paragraph (Chr(13) is PowerPoint's paragraph
#    break character).
$tr.InsertAfter([char]13 + "test")

# 2) Split the original "For test" paragraph into two runs: "For " and
#    "test". Re-assigning the Text of a Characters() sub-range forces the
#    underlying run to split at that boundary.
$para1 = $tr.Paragraphs(1, 1)
$firstWord = $para1.Characters(1, 4)
$firstWord.Text = "For "
